$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the sample data row (Date "2015.11.07" + class/topic/subtopic/
# definition/quiz/exercise/sum counts) - only the header row stays.
$ws.Range("A2:H2").EntireRow.Delete()

# Remove the extra statistics columns (Classes, Topics, Subtopics,
# Definitions, Quiz, Sum) - only Date/Exercises columns remain.
$ws.Range("B1:H1").EntireColumn.Delete()

# Re-create column B (for the "Exercises" header) with default width.
$ws.Columns("B").Insert()
$ws.Range("B1").Value = "Exercises"

# Select the rows below the header, as left by the author (prepping
# the new navbar-menu content).
[void]$ws.Rows("2:3").Select()
